# QE_holdings.xlsx update
#  - Bump the "Model holdings provided as of ..." date in the confidential
#    disclaimer text from 2021-06-14 to 2021-07-07.
#  - Refresh the Weight (D) and Percent Change (E) columns for every holding
#    row (2-34) plus the Total row (35) with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the cells can be written, then
# restore protection afterwards.
$ws.Unprotect()

# --- Disclaimer text (shared string) ---------------------------------
$oldDate = "2021-06-14"
$newDate = "2021-07-07"
$ws.Range("A38").Replace($oldDate, $newDate)

# --- Weight / Percent Change values -----------------------------------
$values = @{
    2  = @(0.0933253620812682,  0)
    3  = @(0.07932967198583461, 0)
    4  = @(0.05250940872023548, 0)
    5  = @(0.05183542364419073, -0.006448839208942458)
    6  = @(0.04791893893778066, 0)
    7  = @(0.0416389351378648,  0)
    8  = @(0.036053487387065,   0)
    9  = @(0.03839229274208928, 0)
    10 = @(0.03398467605360167, 0)
    11 = @(0.03516436421764946, 0)
    12 = @(0.03527721886158235, 0)
    13 = @(0.03091045841127804, 0)
    14 = @(0.03192529308279602, 0)
    15 = @(0.03223814329825554, 0)
    16 = @(0.02977476927296609, 0)
    17 = @(0.02913792623163344, 0)
    18 = @(0.02826708837161463, 0)
    19 = @(0.02388718535518142, 0)
    20 = @(0.02095182178108913, 0)
    21 = @(0.02191737182962377, 0)
    22 = @(0.0215515227877097,  0)
    23 = @(0.02083296727001045, 0)
    24 = @(0.01906986345302354, 0)
    25 = @(0.02136295553455602, 0)
    26 = @(0.0202058382992948,  0)
    27 = @(0.0196778499904645,  0)
    28 = @(0.0186713008497669,  0)
    29 = @(0.02058511561529709, 0)
    30 = @(0.0118274523921256,  0)
    31 = @(0.008360957721650448, 0)
    32 = @(0.007737542954406009, 0)
    33 = @(0.008526954046017552, 0)
    34 = @(0.007149841682077039, 0)
    35 = @(1, -0.0003342783124089532)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("D$row").Value = $pair[0]
    $ws.Range("E$row").Value = $pair[1]
}

# Restore sheet protection.
$ws.Protect()
